# Updates rows 117-120 on the active sheet ("Uruguay Primera División")
# This performs a circular rotation of the 4 match records:
#   new row117 <- old row120 data
#   new row118 <- old row117 data
#   new row119 <- old row118 data
#   new row120 <- old row119 data
# Column A (row index) stays as-is for each row; columns B and E:AD are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for columns B and E:AD across rows 117-120
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")
$rows = @(117,118,119,120)

$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Rotation mapping: destination row -> source row (where the new data comes from)
$rotation = @{ 117 = 120; 118 = 117; 119 = 118; 120 = 119 }

foreach ($dst in $rows) {
    $src = $rotation[$dst]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value = $srcVals[$c]
    }
}
